# Update header labels in row 2 (shared-string backed cells)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Feature Type"
$ws.Range("C2").Value = "Sample Features"

# Widen column B (was 25.01 characters, now ~29.37 characters)
$ws.Columns.Item(2).ColumnWidth = 28.5

# Update the view: zoom in to 140% and change the active selection to B2:C9
$ws.Activate()
[void]$ws.Range("B2:C9").Select()
$excel.ActiveWindow.Zoom = 140
